{"js": "// Replace the date line and each division problem in the practice-sheet\n// table with the values from the new day's worksheet.\nconst replacements = [\n  [\"2025-05-17 Saturday\", \"2025-05-18 Sunday\"],\n  [\"591\u00f74=\", \"555\u00f75=\"],\n  [\"180\u00f74=\", \"404\u00f74=\"],\n  [\"520\u00f74=\", \"589\u00f78=\"],\n  [\"556\u00f78=\", \"953\u00f72=\"],\n  [\"363\u00f76=\", \"733\u00f79=\"],\n  [\"853\u00f77=\", \"591\u00f76=\"],\n  [\"966\u00f72=\", \"861\u00f73=\"],\n  [\"627\u00f77=\", \"980\u00f77=\"],\n  [\"619\u00f73=\", \"217\u00f76=\"],\n  [\"129\u00f74=\", \"476\u00f74=\"],\n  [\"905\u00f77=\", \"781\u00f73=\"],\n  [\"817\u00f76=\", \"853\u00f73=\"],\n  [\"329\u00f74=\", \"834\u00f75=\"],\n  [\"569\u00f79=\", \"508\u00f78=\"],\n  [\"802\u00f78=\", \"180\u00f77=\"],\n  [\"738\u00f72=\", \"103\u00f74=\"],\n  [\"464\u00f78=\", \"770\u00f73=\"],\n  [\"962\u00f72=\", \"731\u00f73=\"],\n  [\"975\u00f74=\", \"327\u00f72=\"],\n  [\"373\u00f76=\", \"573\u00f75=\"],\n  [\"290\u00f72=\", \"181\u00f79=\"],\n  [\"631\u00f79=\", \"389\u00f76=\"],\n  [\"635\u00f76=\", \"554\u00f73=\"],\n  [\"801\u00f77=\", \"907\u00f77=\"],\n  [\"837\u00f75=\", \"572\u00f79=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each division problem in the practice-sheet\n# table with the values from the new day's worksheet.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-05-17 Saturday\", \"2025-05-18 Sunday\"),\n    @(\"591\u00f74=\", \"555\u00f75=\"),\n    @(\"180\u00f74=\", \"404\u00f74=\"),\n    @(\"520\u00f74=\", \"589\u00f78=\"),\n    @(\"556\u00f78=\", \"953\u00f72=\"),\n    @(\"363\u00f76=\", \"733\u00f79=\"),\n    @(\"853\u00f77=\", \"591\u00f76=\"),\n    @(\"966\u00f72=\", \"861\u00f73=\"),\n    @(\"627\u00f77=\", \"980\u00f77=\"),\n    @(\"619\u00f73=\", \"217\u00f76=\"),\n    @(\"129\u00f74=\", \"476\u00f74=\"),\n    @(\"905\u00f77=\", \"781\u00f73=\"),\n    @(\"817\u00f76=\", \"853\u00f73=\"),\n    @(\"329\u00f74=\", \"834\u00f75=\"),\n    @(\"569\u00f79=\", \"508\u00f78=\"),\n    @(\"802\u00f78=\", \"180\u00f77=\"),\n    @(\"738\u00f72=\", \"103\u00f74=\"),\n    @(\"464\u00f78=\", \"770\u00f73=\"),\n    @(\"962\u00f72=\", \"731\u00f73=\"),\n    @(\"975\u00f74=\", \"327\u00f72=\"),\n    @(\"373\u00f76=\", \"573\u00f75=\"),\n    @(\"290\u00f72=\", \"181\u00f79=\"),\n    @(\"631\u00f79=\", \"389\u00f76=\"),\n    @(\"635\u00f76=\", \"554\u00f73=\"),\n    @(\"801\u00f77=\", \"907\u00f77=\"),\n    @(\"837\u00f75=\", \"572\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
